# The sheet holds a small "Rect(...)" source-code generator table:
#   columns L, M, N, O  -> raw source numbers for each row
#   column  J           -> the row's sequence index
#   columns B, D, F, H  -> shared formulas (=M, =L, =M+O, =L+N) that
#                           recompute automatically once L/M/N/O change
#
# This updates the raw values for rows 1-10 to the new numbers, and
# renumbers J8:J10 (the table previously skipped the value 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @{ Row = 1;  L = 0;  M = 75; N = 31; O = 38 }
    @{ Row = 2;  L = 37; M = 0;  N = 32; O = 39 }
    @{ Row = 3;  L = 0;  M = 33; N = 31; O = 41 }
    @{ Row = 4;  L = 70; M = 0;  N = 31; O = 42 }
    @{ Row = 5;  L = 32; M = 40; N = 30; O = 38 }
    @{ Row = 6;  L = 0;  M = 0;  N = 36; O = 32 }
    @{ Row = 7;  L = 36; M = 68; N = 33; O = 35 }
    @{ Row = 8;  L = 0;  M = 34; N = 35; O = 34; J = 7 }
    @{ Row = 9;  L = 39; M = 0;  N = 37; O = 33; J = 8 }
    @{ Row = 10; L = 0;  M = 69; N = 35; O = 33; J = 9 }
)

foreach ($entry in $newValues) {
    $r = $entry.Row
    $ws.Cells.Item($r, 12).Value = $entry.L   # column L
    $ws.Cells.Item($r, 13).Value = $entry.M   # column M
    $ws.Cells.Item($r, 14).Value = $entry.N   # column N
    $ws.Cells.Item($r, 15).Value = $entry.O   # column O
    if ($entry.ContainsKey("J")) {
        $ws.Cells.Item($r, 10).Value = $entry.J   # column J
    }
}

[void]$excel.Calculate()

# Reflect the extended selection recorded in the sheet view (A1:I2 -> A1:I6).
[void]$ws.Range("A1:I6").Select()
